$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from E1 (geneConfidence/panel header) onto F1 so the
# new "time_taken" header matches the existing bold/bordered/centered style.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate the per-row "time_taken" timestamps for rows 2-13.
$timestamps = @(
    "2021-10-05 10:52:41.899352",
    "2021-10-05 10:52:41.899363",
    "2021-10-05 10:52:41.899367",
    "2021-10-05 10:52:41.899370",
    "2021-10-05 10:52:41.899374",
    "2021-10-05 10:52:41.899377",
    "2021-10-05 10:52:41.899380",
    "2021-10-05 10:52:41.899383",
    "2021-10-05 10:52:41.899386",
    "2021-10-05 10:52:41.899389",
    "2021-10-05 10:52:41.899392",
    "2021-10-05 10:52:41.899395"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
